$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Collected" column (D) entirely - values were all a flat ₹0.00
# placeholder and are no longer needed. Remaining columns (SPV, Phone) shift
# left from E:F into D:E.
$ws.Columns("D").Delete() | Out-Null

# Leave the entire (now-empty-of-special-meaning) column D selected, matching
# the state Excel leaves behind right after a column delete.
$ws.Columns("D").Select() | Out-Null

